# Template_GHG_GCHeadspace.xlsx -- "Added columns for air averages to template"
#
# 1) Insert two new header columns ("Site", "Sample_Type") immediately to the
#    left of the existing "Site_ID" column (old column I), pushing every
#    column from the old I onward two slots to the right.
# 2) Both new header cells carry the placeholder note "Need to add in old
#    samplings" in row 2 (mirroring the existing WaterT_C/Air_Location note
#    pattern already on the sheet).
# 3) Append six brand-new trailing columns capturing air CO2/CH4 summary
#    stats (med/min/max) right after the current last column (WaterV_mL).
# 4) The threaded comment that used to sit on "Site_ID" (old I1) has to stay
#    attached to that same header text, which is now at K1 -- column
#    insertion shifts cell contents but Excel does NOT relocate existing
#    comments, so we move it explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert "Site" / "Sample_Type" ahead of Site_ID -----------------
$ws.Columns("I:J").Insert()

$ws.Range("I1").Value = "Site"
$ws.Range("J1").Value = "Sample_Type"

# --- 2) Row-2 placeholder note on both new headers ----------------------
$ws.Range("I2").Value = "Need to add in old samplings"
$ws.Range("J2").Value = "Need to add in old samplings"

# --- 3) New trailing columns for air CO2/CH4 summary stats --------------
$ws.Range("Y1").Value = "AirCO2_med_ppm"
$ws.Range("Z1").Value = "AirCO2_min_ppm"
$ws.Range("AA1").Value = "AirCO2_max_ppm"
$ws.Range("AB1").Value = "AirCH4_med_ppm"
$ws.Range("AC1").Value = "AirCH4_min_ppm"
$ws.Range("AD1").Value = "AirCH4_max_ppm"

# --- 4) Re-home the threaded comment from old I1 (Site_ID) to new K1 ----
$movedComment = $ws.Range("I1").Comment
if ($movedComment -ne $null) {
    $commentText = $movedComment.Text()
    $movedComment.Delete()
    $ws.Range("K1").AddCommentThreaded($commentText)
}

# --- Cosmetic touch-ups mirroring the author's on-screen state ----------
# (best effort -- exact pixel-perfect fit metrics from the original
# authoring machine can't be reproduced, but approximate the same visual
# result: a taller wrapped header row and reasonably sized new columns)
$ws.Rows("1").RowHeight = 43.2
$ws.Columns("I").ColumnWidth = 24.5546875
$ws.Columns("J").ColumnWidth = 12.33203125

$ws.Range("U9").Select()
